# fix(gui) step 1 and 2
# - Bump the date in A1 by one day (45308 -> 45309)
# - Update the unit prices in D30:D37 (steps 1 & 2 price list)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 45309

$ws.Range("D30").Value = 120.069
$ws.Range("D31").Value = 128.588
$ws.Range("D32").Value = 140.928
$ws.Range("D33").Value = 208.099
$ws.Range("D34").Value = 269.256
$ws.Range("D35").Value = 379.682
$ws.Range("D36").Value = 550.922
$ws.Range("D37").Value = 844.688
